$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.713252999999999
$ws.Range("H2").Value = 20.139759
$ws.Range("I2").Value = 0.3101840064655811
$ws.Range("J2").Value = 0.3231642354899327
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.474639
$ws.Range("N2").Value = 4.423916999999999
$ws.Range("O2").Value = 0.02436078515116209
$ws.Range("P2").Value = 0.02448628780347778
$ws.Range("Q2").Value = 9.899624690666997
$ws.Range("R2").Value = 89.09662221600298
$ws.Range("S2").Value = 0.007556325938834693
$ws.Range("T2").Value = 0.007913092477997359
$ws.Range("G3").Value = 6.713252999999999
$ws.Range("H3").Value = 20.139759
$ws.Range("I3").Value = 0.3101840064655811
$ws.Range("J3").Value = 0.3231642354899327
$ws.Range("O3").Value = 0.01812262212538788
$ws.Range("P3").Value = 0.01821598681497149
$ws.Range("Q3").Value = 7.364588470316999
$ws.Range("R3").Value = 66.28129623285299
$ws.Range("S3").Value = 0.005621347538514596
$ws.Range("T3").Value = 0.005886755452754954
$ws.Range("G4").Value = 6.713252999999999
$ws.Range("H4").Value = 20.139759
$ws.Range("I4").Value = 0.3101840064655811
$ws.Range("J4").Value = 0.3231642354899327
$ws.Range("M4").Value = 36.94558466666667
$ws.Range("N4").Value = 110.836754
$ws.Range("O4").Value = 0.6103347669149772
$ws.Range("P4").Value = 0.6134791085925136
$ws.Range("Q4").Value = 248.025057100254
$ws.Range("R4").Value = 2232.225513902286
$ws.Range("S4").Value = 0.1893160832869243
$ws.Range("T4").Value = 0.198254507117345
$ws.Range("G5").Value = 6.713252999999999
$ws.Range("H5").Value = 20.139759
$ws.Range("I5").Value = 0.3101840064655811
$ws.Range("J5").Value = 0.3231642354899327
$ws.Range("M5").Value = 0.930777
$ws.Range("N5").Value = 1.861554
$ws.Range("O5").Value = 0.01537627753005528
$ws.Range("P5").Value = 0.01030366234396244
$ws.Range("Q5").Value = 6.248541487580999
$ws.Range("R5").Value = 37.491248925486
$ws.Range("S5").Value = 0.004769475368799238
$ws.Range("T5").Value = 0.003329775164133029
$ws.Range("G6").Value = 6.713252999999999
$ws.Range("H6").Value = 20.139759
$ws.Range("I6").Value = 0.3101840064655811
$ws.Range("J6").Value = 0.3231642354899327
$ws.Range("M6").Value = 20.08528866666667
$ws.Range("N6").Value = 60.255866
$ws.Range("O6").Value = 0.3318055482784176
$ws.Range("P6").Value = 0.3335149544450747
$ws.Range("Q6").Value = 134.837624397366
$ws.Range("R6").Value = 1213.538619576294
$ws.Range("S6").Value = 0.1029207743325084
$ws.Range("T6").Value = 0.1077801052777023
$ws.Range("I7").Value = 0.0154484264788496
$ws.Range("J7").Value = 0.01609489473505086
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.474639
$ws.Range("N7").Value = 4.423916999999999
$ws.Range("O7").Value = 0.02436078515116209
$ws.Range("P7").Value = 0.02448628780347778
$ws.Range("Q7").Value = 0.4930416172793333
$ws.Range("R7").Value = 4.437374555513999
$ws.Range("S7").Value = 0.0003763357983747784
$ws.Range("T7").Value = 0.0003941042246491345
$ws.Range("I8").Value = 0.0154484264788496
$ws.Range("J8").Value = 0.01609489473505086
$ws.Range("O8").Value = 0.01812262212538788
$ws.Range("P8").Value = 0.01821598681497149
$ws.Range("S8").Value = 0.0002799659955080276
$ws.Range("T8").Value = 0.0002931843902820404
$ws.Range("I9").Value = 0.0154484264788496
$ws.Range("J9").Value = 0.01609489473505086
$ws.Range("M9").Value = 36.94558466666667
$ws.Range("N9").Value = 110.836754
$ws.Range("O9").Value = 0.6103347669149772
$ws.Range("P9").Value = 0.6134791085925136
$ws.Range("Q9").Value = 12.35265771174089
$ws.Range("R9").Value = 111.173919405668
$ws.Range("S9").Value = 0.00942871177417183
$ws.Range("T9").Value = 0.00987388167494934
$ws.Range("I10").Value = 0.0154484264788496
$ws.Range("J10").Value = 0.01609489473505086
$ws.Range("M10").Value = 0.930777
$ws.Range("N10").Value = 1.861554
$ws.Range("O10").Value = 0.01537627753005528
$ws.Range("P10").Value = 0.01030366234396244
$ws.Range("Q10").Value = 0.311202807878
$ws.Range("R10").Value = 1.867216847268
$ws.Range("S10").Value = 0.0002375392929414461
$ws.Range("T10").Value = 0.0001658363608115828
$ws.Range("I11").Value = 0.0154484264788496
$ws.Range("J11").Value = 0.01609489473505086
$ws.Range("M11").Value = 20.08528866666667
$ws.Range("N11").Value = 60.255866
$ws.Range("O11").Value = 0.3318055482784176
$ws.Range("P11").Value = 0.3335149544450747
$ws.Range("Q11").Value = 6.715462704930222
$ws.Range("R11").Value = 60.439164344372
$ws.Range("S11").Value = 0.005125873617853515
$ws.Range("T11").Value = 0.005367888084358758
$ws.Range("G12").Value = 6.661784666666667
$ws.Range("H12").Value = 19.985354
$ws.Range("I12").Value = 0.3078059262949933
$ws.Range("J12").Value = 0.3206866401135023
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 1.474639
$ws.Range("N12").Value = 4.423916999999999
$ws.Range("O12").Value = 0.02436078515116209
$ws.Range("P12").Value = 0.02448628780347778
$ws.Range("Q12").Value = 9.823727479068666
$ws.Range("R12").Value = 88.41354731161799
$ws.Range("S12").Value = 0.007498394038726764
$ws.Range("T12").Value = 0.007852425364549519
$ws.Range("G13").Value = 6.661784666666667
$ws.Range("H13").Value = 19.985354
$ws.Range("I13").Value = 0.3078059262949933
$ws.Range("J13").Value = 0.3206866401135023
$ws.Range("O13").Value = 0.01812262212538788
$ws.Range("P13").Value = 0.01821598681497149
$ws.Range("Q13").Value = 7.308126559190889
$ws.Range("R13").Value = 65.77313903271801
$ws.Range("S13").Value = 0.005578250490199155
$ws.Range("T13").Value = 0.005841623608045064
$ws.Range("G14").Value = 6.661784666666667
$ws.Range("H14").Value = 19.985354
$ws.Range("I14").Value = 0.3078059262949933
$ws.Range("J14").Value = 0.3206866401135023
$ws.Range("M14").Value = 36.94558466666667
$ws.Range("N14").Value = 110.836754
$ws.Range("O14").Value = 0.6103347669149772
$ws.Range("P14").Value = 0.6134791085925136
$ws.Range("Q14").Value = 246.1235294334351
$ws.Range("R14").Value = 2215.111764900916
$ws.Range("S14").Value = 0.1878646582803034
$ws.Range("T14").Value = 0.1967345541143596
$ws.Range("G15").Value = 6.661784666666667
$ws.Range("H15").Value = 19.985354
$ws.Range("I15").Value = 0.3078059262949933
$ws.Range("J15").Value = 0.3206866401135023
$ws.Range("M15").Value = 0.930777
$ws.Range("N15").Value = 1.861554
$ws.Range("O15").Value = 0.01537627753005528
$ws.Range("P15").Value = 0.01030366234396244
$ws.Range("Q15").Value = 6.200635946686
$ws.Range("R15").Value = 37.203815680116
$ws.Range("S15").Value = 0.004732909348107559
$ws.Range("T15").Value = 0.003304246857949328
$ws.Range("G16").Value = 6.661784666666667
$ws.Range("H16").Value = 19.985354
$ws.Range("I16").Value = 0.3078059262949933
$ws.Range("J16").Value = 0.3206866401135023
$ws.Range("M16").Value = 20.08528866666667
$ws.Range("N16").Value = 60.255866
$ws.Range("O16").Value = 0.3318055482784176
$ws.Range("P16").Value = 0.3335149544450747
$ws.Range("Q16").Value = 133.8038680651738
$ws.Range("R16").Value = 1204.234812586564
$ws.Range("S16").Value = 0.1021317141376565
$ws.Range("T16").Value = 0.1069537901685988
$ws.Range("G17").Value = 2.607918
$ws.Range("H17").Value = 5.215835999999999
$ws.Range("I17").Value = 0.1204981331366039
$ws.Range("J17").Value = 0.08369373503331734
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 1.474639
$ws.Range("N17").Value = 4.423916999999999
$ws.Range("O17").Value = 0.02436078515116209
$ws.Range("P17").Value = 0.02448628780347778
$ws.Range("Q17").Value = 3.845737591601999
$ws.Range("R17").Value = 23.07442554961199
$ws.Range("S17").Value = 0.002935429132456932
$ws.Range("T17").Value = 0.002049348883373819
$ws.Range("G18").Value = 2.607918
$ws.Range("H18").Value = 5.215835999999999
$ws.Range("I18").Value = 0.1204981331366039
$ws.Range("J18").Value = 0.08369373503331734
$ws.Range("O18").Value = 0.01812262212538788
$ws.Range("P18").Value = 0.01821598681497149
$ws.Range("Q18").Value = 2.860944289502
$ws.Range("R18").Value = 17.165665737012
$ws.Range("S18").Value = 0.002183742133649351
$ws.Range("T18").Value = 0.001524563973862626
$ws.Range("G19").Value = 2.607918
$ws.Range("H19").Value = 5.215835999999999
$ws.Range("I19").Value = 0.1204981331366039
$ws.Range("J19").Value = 0.08369373503331734
$ws.Range("M19").Value = 36.94558466666667
$ws.Range("N19").Value = 110.836754
$ws.Range("O19").Value = 0.6103347669149772
$ws.Range("P19").Value = 0.6134791085925136
$ws.Range("Q19").Value = 96.351055272724
$ws.Range("R19").Value = 578.106331636344
$ws.Range("S19").Value = 0.07354420000161903
$ws.Range("T19").Value = 0.05134435796301755
$ws.Range("G20").Value = 2.607918
$ws.Range("H20").Value = 5.215835999999999
$ws.Range("I20").Value = 0.1204981331366039
$ws.Range("J20").Value = 0.08369373503331734
$ws.Range("M20").Value = 0.930777
$ws.Range("N20").Value = 1.861554
$ws.Range("O20").Value = 0.01537627753005528
$ws.Range("P20").Value = 0.01030366234396244
$ws.Range("Q20").Value = 2.427390092286
$ws.Range("R20").Value = 9.709560369143999
$ws.Range("S20").Value = 0.001852812736961972
$ws.Range("T20").Value = 0.0008623519860883618
$ws.Range("G21").Value = 2.607918
$ws.Range("H21").Value = 5.215835999999999
$ws.Range("I21").Value = 0.1204981331366039
$ws.Range("J21").Value = 0.08369373503331734
$ws.Range("M21").Value = 20.08528866666667
$ws.Range("N21").Value = 60.255866
$ws.Range("O21").Value = 0.3318055482784176
$ws.Range("P21").Value = 0.3335149544450747
$ws.Range("Q21").Value = 52.380785848996
$ws.Range("R21").Value = 314.2847150939759
$ws.Range("S21").Value = 0.03998194913191662
$ws.Range("T21").Value = 0.02791311222697498
$ws.Range("G22").Value = 5.325505333333333
$ws.Range("H22").Value = 15.976516
$ws.Range("I22").Value = 0.2460635076239721
$ws.Range("J22").Value = 0.2563604946281968
$ws.Range("K22").Value = 3
$ws.Range("L22").Value = 1
$ws.Range("M22").Value = 1.474639
$ws.Range("N22").Value = 4.423916999999999
$ws.Range("O22").Value = 0.02436078515116209
$ws.Range("P22").Value = 0.02448628780347778
$ws.Range("Q22").Value = 7.853197859241332
$ws.Range("R22").Value = 70.67878073317199
$ws.Range("S22").Value = 0.005994300242768919
$ws.Range("T22").Value = 0.006277316852907945
$ws.Range("G23").Value = 5.325505333333333
$ws.Range("H23").Value = 15.976516
$ws.Range("I23").Value = 0.2460635076239721
$ws.Range("J23").Value = 0.2563604946281968
$ws.Range("O23").Value = 0.01812262212538788
$ws.Range("P23").Value = 0.01821598681497149
$ws.Range("Q23").Value = 5.842198286952444
$ws.Range("R23").Value = 52.579784582572
$ws.Range("S23").Value = 0.004459315967516745
$ws.Range("T23").Value = 0.004669859390026801
$ws.Range("G24").Value = 5.325505333333333
$ws.Range("H24").Value = 15.976516
$ws.Range("I24").Value = 0.2460635076239721
$ws.Range("J24").Value = 0.2563604946281968
$ws.Range("M24").Value = 36.94558466666667
$ws.Range("N24").Value = 110.836754
$ws.Range("O24").Value = 0.6103347669149772
$ws.Range("P24").Value = 0.6134791085925136
$ws.Range("Q24").Value = 196.7539081854516
$ws.Range("R24").Value = 1770.785173669064
$ws.Range("S24").Value = 0.1501811135719587
$ws.Range("T24").Value = 0.157271807722842
$ws.Range("G25").Value = 5.325505333333333
$ws.Range("H25").Value = 15.976516
$ws.Range("I25").Value = 0.2460635076239721
$ws.Range("J25").Value = 0.2563604946281968
$ws.Range("M25").Value = 0.930777
$ws.Range("N25").Value = 1.861554
$ws.Range("O25").Value = 0.01537627753005528
$ws.Range("P25").Value = 0.01030366234396244
$ws.Range("Q25").Value = 4.956857877644
$ws.Range("R25").Value = 29.741147265864
$ws.Range("S25").Value = 0.003783540783245069
$ws.Range("T25").Value = 0.002641451974980136
$ws.Range("G26").Value = 5.325505333333333
$ws.Range("H26").Value = 15.976516
$ws.Range("I26").Value = 0.2460635076239721
$ws.Range("J26").Value = 0.2563604946281968
$ws.Range("M26").Value = 20.08528866666667
$ws.Range("N26").Value = 60.255866
$ws.Range("O26").Value = 0.3318055482784176
$ws.Range("P26").Value = 0.3335149544450747
$ws.Range("Q26").Value = 106.9643119158729
$ws.Range("R26").Value = 962.678807242856
$ws.Range("S26").Value = 0.08164523705848266
$ws.Range("T26").Value = 0.08550005868743986
